# Update "想去人数" (F column) counts on the "展览", "演出" and "全部类型" sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3179
$ws1.Range("F3").Value = 733
$ws1.Range("F4").Value = 118
$ws1.Range("F5").Value = 6866
$ws1.Range("F6").Value = 1989
$ws1.Range("F7").Value = 21
$ws1.Range("F8").Value = 71
$ws1.Range("F12").Value = 19
$ws1.Range("F14").Value = 176
$ws1.Range("F15").Value = 34

# --- Sheet "演出" (Performance) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 15

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3179
$ws4.Range("F3").Value = 15
$ws4.Range("F4").Value = 733
$ws4.Range("F5").Value = 118
$ws4.Range("F6").Value = 6866
$ws4.Range("F7").Value = 1989
$ws4.Range("F8").Value = 21
$ws4.Range("F9").Value = 71
$ws4.Range("F13").Value = 19
$ws4.Range("F15").Value = 176
$ws4.Range("F16").Value = 34
